$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the report template placeholder: {{Items.LectureName}} -> {{Items.LecturerName}}
$ws.Range("E3").Value = "{{Items.LecturerName}}"

# Move the active selection to E3 (matches the saved sheet view state)
$ws.Range("E3").Select()
